$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '39.873.29'
$ws.Range('E2').Value = '  +0.88%  '
$ws.Range('D3').Value = '2.244.64'
$ws.Range('E3').Value = '  -2.89%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '298.79'
$ws.Range('E5').Value = '  -2.17%  '
$ws.Range('D6').Value = '84.23'
$ws.Range('E6').Value = '  +0.89%  '
$ws.Range('D7').Value = '0.518'
$ws.Range('E7').Value = '  -1.67%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').Value = '0.477'
$ws.Range('E9').Value = '  -0.17%  '
$ws.Range('D10').Value = '30.50'
$ws.Range('E10').Value = '  +3.48%  '
$ws.Range('E11').Value = '  -2.70%  '
$ws.Range('D12').Value = '47.19'
$ws.Range('E12').Value = '  -9.79%  '
$ws.Range('E13').Value = '  -1.84%  '
$ws.Range('D14').Value = '2.588.22'
$ws.Range('E14').Value = '  -3.09%  '
$ws.Range('D15').Value = '6.36'
$ws.Range('E15').Value = '  +0.38%  '
$ws.Range('D16').Value = '14.28'
$ws.Range('E16').Value = '  -1.85%  '
$ws.Range('D17').Value = '2.231.23'
$ws.Range('E17').Value = '  -3.33%  '
$ws.Range('D18').Value = '0.724'
$ws.Range('E18').Value = '  -2.90%  '
$ws.Range('D19').Value = '39.793.56'
$ws.Range('E19').Value = '  +0.65%  '
$ws.Range('D20').Value = '0.0₃0883'
$ws.Range('E20').Value = '  -0.91%  '
$ws.Range('D21').Value = '5.84'
$ws.Range('E21').Value = '  -3.31%  '
$ws.Range('D22').Value = '65.49'
$ws.Range('E22').Value = '  -2.68%  '
$ws.Range('D23').Value = '10.49'
$ws.Range('E23').Value = '  +0.26%  '
$ws.Range('D24').Value = '229.69'
$ws.Range('E24').Value = '  -2.14%  '
$ws.Range('E25').Value = '  -0.09%  '
$ws.Range('B26').Value = 'ImmutableX'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D26').Value = '1.87'
$ws.Range('E26').Value = '  +5.29%  '
$ws.Range('B27').Value = 'LEO'
$ws.Range('C27').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D27').Value = '4.01'
$ws.Range('E27').Value = '  +1.80%  '
$ws.Range('B28').Value = 'PancakeSwap'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D28').Value = '2.43'
$ws.Range('E28').Value = '  -3.57%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = '23.10'
$ws.Range('E29').Value = '  +0.87%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = '2.18'
$ws.Range('E30').Value = '  +3.07%  '
$ws.Range('B31').Value = 'Cosmos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D31').Value = '9.22'
$ws.Range('E31').Value = '  +0.51%  '
$ws.Range('B32').Value = 'InjectiveProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D32').Value = '32.94'
$ws.Range('E32').Value = '  -2.57%  '
$ws.Range('B33').Value = 'Monero'
$ws.Range('C33').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D33').Value = '150.13'
$ws.Range('E33').Value = '  +0.11%  '
$ws.Range('B34').Value = 'FirstDigitalUSD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  -0.31%  '
$ws.Range('E35').Value = '  -2.86%  '
$ws.Range('B36').Value = 'WEMIXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').Value = '2.43'
$ws.Range('E36').Value = '  -0.57%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').Value = '0.0707'
$ws.Range('E37').Value = '  -0.60%  '
$ws.Range('B38').Value = 'Celestia'
$ws.Range('C38').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D38').Value = '16.39'
$ws.Range('E38').Value = '  +7.61%  '
$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D39').Value = '0.112'
$ws.Range('E39').Value = '  -0.90%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').Value = '0.0979'
$ws.Range('E40').Value = '  -0.17%  '
$ws.Range('B41').Value = 'LidoDAOToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D41').Value = '2.69'
$ws.Range('E41').Value = '  -0.48%  '
$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').Value = '1.69'
$ws.Range('E42').Value = '  +0.48%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').Value = '3.75'
$ws.Range('E43').Value = '  +0.09%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '1.930.50'
$ws.Range('E44').Value = '  +0.03%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').Value = '0.0266'
$ws.Range('E45').Value = '  +1.84%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').Value = '2.05'
$ws.Range('E46').Value = '  -8.73%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '16.57'
$ws.Range('E47').Value = '  -4.71%  '
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D48').Value = '9.16'
$ws.Range('E48').Value = '  -1.27%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').Value = '2.65'
$ws.Range('E49').Value = '  +0.20%  '
$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').Value = '2.457.00'
$ws.Range('E50').Value = '  -3.18%  '
$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D51').Value = '72.15'
